$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Stamp number-format styles onto the new rows (3-6) by copying
#    formats from existing styled cells, so we reuse existing style
#    indices (date / text-doc-number / text-amount) instead of the
#    COM layer minting brand-new style entries for every new cell.
# ------------------------------------------------------------------
$ws.Cells.Item(1,1).Copy() | Out-Null
$ws.Range("A3:A6").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,4).Copy() | Out-Null
$ws.Range("D3:D6").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(2,6).Copy() | Out-Null
$ws.Cells.Item(6,6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(2,7).Copy() | Out-Null
$ws.Cells.Item(6,7).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Pre-mark the numeric-looking amount cells that must stay plain
#    (default-styled) TEXT as "@" so assigning "180.00" etc. does not
#    get auto-coerced into the number 180. Cleared back to the Normal
#    style afterwards so no stray style index is left on the cell.
#    (F6/G6 already get a real text style from the paste above, so
#    they are excluded here.)
# ------------------------------------------------------------------
foreach ($addr in @("F1","G1","F2","G2","F3","G3","F4","G4","F5","G5")) {
  $ws.Range($addr).NumberFormat = "@"
}

# --- Rows 1-6, columns A-G ---
# Row 1
$ws.Cells.Item(1,1).Value = 41822
$ws.Cells.Item(1,2).Value = "RETIROS ATM CLIENTES PRODUBANCO"
$ws.Cells.Item(1,3).Value = "D"
$ws.Cells.Item(1,4).Value = "9382000143"
$ws.Cells.Item(1,5).Value = "AG. PDBCO EXPRESS MEGAMAXI UIO"
$ws.Cells.Item(1,6).Value = "180.00"
$ws.Cells.Item(1,7).Value = "4.30"

# Row 2
$ws.Cells.Item(2,1).Value = 41822
$ws.Cells.Item(2,2).Value = "RETIRO CON LIBRETA"
$ws.Cells.Item(2,3).Value = "D"
$ws.Cells.Item(2,4).Value = "0028525948"
$ws.Cells.Item(2,5).Value = "AG. PDBCO EXPRESS MEGAMAXI UIO"
$ws.Cells.Item(2,6).Value = "1000.00"
$ws.Cells.Item(2,7).Value = "184.30"

# Row 3
$ws.Cells.Item(3,1).Value = 41820
$ws.Cells.Item(3,2).Value = "CAPITALIZACION DE INTERESES EN CUENTA"
$ws.Cells.Item(3,3).Value = "C"
$ws.Cells.Item(3,4).Value = "0020140630"
$ws.Cells.Item(3,5).Value = "MATRIZ - QUITO"
$ws.Cells.Item(3,6).Value = "0.08"
$ws.Cells.Item(3,7).Value = "1184.30"

# Row 4
$ws.Cells.Item(4,1).Value = 41820
$ws.Cells.Item(4,2).Value = "ROL DE PAGOS"
$ws.Cells.Item(4,3).Value = "C"
$ws.Cells.Item(4,4).Value = "0199880577"
$ws.Cells.Item(4,5).Value = "MATRIZ - QUITO"
$ws.Cells.Item(4,6).Value = "1175.18"
$ws.Cells.Item(4,7).Value = "1184.22"

# Row 5
$ws.Cells.Item(5,1).Value = 41802
$ws.Cells.Item(5,2).Value = "RETIRO CON LIBRETA"
$ws.Cells.Item(5,3).Value = "D"
$ws.Cells.Item(5,4).Value = "0031789350"
$ws.Cells.Item(5,5).Value = "AG. PDBCO EXPRESS EL GIRÓN QUI"
$ws.Cells.Item(5,6).Value = "9.04"
$ws.Cells.Item(5,7).Value = "9.04"

# Row 6
$ws.Cells.Item(6,1).Value = 41789
$ws.Cells.Item(6,2).Value = "ROL DE PAGOS"
$ws.Cells.Item(6,3).Value = "C"
$ws.Cells.Item(6,4).Value = "0197325546"
$ws.Cells.Item(6,5).Value = "MATRIZ - QUITO"
$ws.Cells.Item(6,6).Value = "509.04"
$ws.Cells.Item(6,7).Value = "509.04"

# ------------------------------------------------------------------
# 3) Clear the temporary "@" format back to the Normal style on the
#    cells that should end up with no explicit style (s=0), now that
#    the text values are safely stored.
# ------------------------------------------------------------------
foreach ($addr in @("F1","G1","F2","G2","F3","G3","F4","G4","F5","G5")) {
  $ws.Range($addr).Style = "Normal"
}

# --- Column H formulas ---
$ws.Cells.Item(1,8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B1,""', 'mo_tipo' => '"",C1,""', 'mo_documento' => '"",D1,""', 'mo_oficina' => '"",E1,""', 'mo_monto' => "",F1,"", 'mo_saldo' => "",G1,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"
$ws.Cells.Item(2,8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A2,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B2,""', 'mo_tipo' => '"",C2,""', 'mo_documento' => '"",D2,""', 'mo_oficina' => '"",E2,""', 'mo_monto' => "",F2,"", 'mo_saldo' => "",G2,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"
$ws.Cells.Item(3,8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A3,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B3,""', 'mo_tipo' => '"",C3,""', 'mo_documento' => '"",D3,""', 'mo_oficina' => '"",E3,""', 'mo_monto' => "",F3,"", 'mo_saldo' => "",G3,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"
$ws.Cells.Item(4,8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A4,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B4,""', 'mo_tipo' => '"",C4,""', 'mo_documento' => '"",D4,""', 'mo_oficina' => '"",E4,""', 'mo_monto' => "",F4,"", 'mo_saldo' => "",G4,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"
$ws.Cells.Item(5,8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A5,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B5,""', 'mo_tipo' => '"",C5,""', 'mo_documento' => '"",D5,""', 'mo_oficina' => '"",E5,""', 'mo_monto' => "",F5,"", 'mo_saldo' => "",G5,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"
$ws.Cells.Item(6,8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A6,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B6,""', 'mo_tipo' => '"",C6,""', 'mo_documento' => '"",D6,""', 'mo_oficina' => '"",E6,""', 'mo_monto' => "",F6,"", 'mo_saldo' => "",G6,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_fecha_borrado' => NULL, 'mo_quien_borra' => NULL, 'mo_borrado_logico' => false),"")"

# --- View state: selection ---
$ws.Range("H1:H6").Select()
